# Update "想去人数" (interested-count) figures in both the "展览" sheet
# and the consolidated "全部类型" sheet to match the latest scrape.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 540
$wsExhibit.Range("F5").Value = 276
$wsExhibit.Range("F6").Value = 392
$wsExhibit.Range("F7").Value = 241
$wsExhibit.Range("F8").Value = 2305
$wsExhibit.Range("F10").Value = 5731
$wsExhibit.Range("F12").Value = 375

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 540
$wsAll.Range("F6").Value = 276
$wsAll.Range("F7").Value = 392
$wsAll.Range("F8").Value = 241
$wsAll.Range("F11").Value = 2305
$wsAll.Range("F13").Value = 5731
$wsAll.Range("F15").Value = 375
